$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("B2").Value = $false
$ws.Range("C2").Value = 1
$ws.Range("E2").Value = 0.68
$ws.Range("F2").Value = 7

# Add new row 3 with style copied from A2 (style index 1)
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A3").Value = 1

$ws.Range("B3").Value = $true
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.65
$ws.Range("F3").Value = 1
